$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.751.45"
$ws.Range("E2").Value = "  +5.18%  "

# Row 3
$ws.Range("D3").Value = "2.270.26"
$ws.Range("E3").Value = "  +3.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").Value = "'233.39"
$ws.Range("E5").Value = "  +1.79%  "

# Row 6
$ws.Range("D6").Value = "'0.633"
$ws.Range("E6").Value = "  +2.43%  "

# Row 7
$ws.Range("D7").Value = "'63.50"
$ws.Range("E7").Value = "  +5.92%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("E9").Value = "  +7.81%  "

# Row 10
$ws.Range("D10").Value = "'0.104"
$ws.Range("E10").Value = "  +16.87%  "

# Row 11
$ws.Range("D11").Value = "'57.60"
$ws.Range("E11").Value = "  +0.32%  "

# Row 12
$ws.Range("D12").Value = "'26.04"
$ws.Range("E12").Value = "  +17.21%  "

# Row 13
$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "  -0.25%  "

# Row 14
$ws.Range("D14").Value = "2.605.29"
$ws.Range("E14").Value = "  +3.09%  "

# Row 15
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("D16").Value = "'5.96"
$ws.Range("E16").Value = "  +5.83%  "

# Row 17
$ws.Range("D17").Value = "'0.826"
$ws.Range("E17").Value = "  +4.31%  "

# Row 18
$ws.Range("D18").Value = "2.288.52"
$ws.Range("E18").Value = "  +2.72%  "

# Row 19
$ws.Range("D19").Value = "43.605.43"
$ws.Range("E19").Value = "  +5.07%  "

# Row 20
$ws.Range("D20").Value = "'0.0000100"
$ws.Range("E20").Value = "  +11.49%  "

# Row 21
$ws.Range("D21").Value = "'73.99"

# Row 22
$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  +0.84%  "

# Row 23
$ws.Range("D23").Value = "'249.57"
$ws.Range("E23").Value = "  +2.92%  "

# Row 24
$ws.Range("E24").Value = "  +0.12%  "

# Row 25
$ws.Range("D25").Value = "'2.52"
$ws.Range("E25").Value = "  +7.03%  "

# Row 26
$ws.Range("E26").Value = "  -1.69%  "

# Row 27
$ws.Range("E27").Value = "  +2.71%  "

# Row 28
$ws.Range("D28").Value = "'172.76"
$ws.Range("E28").Value = "  +2.38%  "

# Row 29
$ws.Range("E29").Value = "  +6.54%  "

# Row 30
$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = "  -0.62%  "

# Row 31
$ws.Range("E31").Value = "  +2.21%  "

# Row 32
$ws.Range("D32").Value = "'2.80"
$ws.Range("E32").Value = "  +11.13%  "

# Row 33
$ws.Range("E33").Value = "  +1.93%  "

# Row 34
$ws.Range("D34").Value = "'0.0689"
$ws.Range("E34").Value = "  +6.76%  "

# Row 35
$ws.Range("E35").Value = "  +2.06%  "

# Row 36
$ws.Range("D36").Value = "'4.75"
$ws.Range("E36").Value = "  +2.89%  "

# Row 37
$ws.Range("D37").Value = "'6.84"
$ws.Range("E37").Value = "  +6.20%  "

# Row 38
$ws.Range("D38").Value = "'3.84"
$ws.Range("E38").Value = "  +7.87%  "

# Row 39
$ws.Range("E39").Value = "  -1.19%  "

# Row 40
$ws.Range("E40").Value = "  +5.55%  "

# Row 42
$ws.Range("E42").Value = "  -0.94%  "

# Row 43
$ws.Range("D43").Value = "'17.52"
$ws.Range("E43").Value = "  +7.14%  "

# Row 44
$ws.Range("D44").Value = "'10.48"
$ws.Range("E44").Value = "  +22.17%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0962"
$ws.Range("E45").Value = "  +0.85%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.45"
$ws.Range("E46").Value = "  +2.81%  "

# Row 47
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "'1.21"
$ws.Range("E47").Value = "  +0.76%  "

# Row 48
$ws.Range("D48").Value = "'97.80"
$ws.Range("E48").Value = "  +0.69%  "

# Row 49
$ws.Range("D49").Value = "1.480.04"
$ws.Range("E49").Value = "  +1.08%  "

# Row 50
$ws.Range("E50").Value = "  +4.83%  "

# Row 51
$ws.Range("E51").Value = "  +1.22%  "
